# Adds 2023-11-01 crime-count data: bumps the 2023 (column J) year-to-date
# totals across the Citywide Totals, By Neighborhood, and each neighborhood
# detail worksheet, per the source commit "Add data for 2023-11-01".
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 6423
$ws.Range("J3").Value = 6842
$ws.Range("J4").Value = 1478
$ws.Range("J5").Value = 526
$ws.Range("J6").Value = 9048
$ws.Range("J7").Value = 24317

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("J2").Value = 194
$ws.Range("J4").Value = 110
$ws.Range("J6").Value = 185
$ws.Range("J7").Value = 708
$ws.Range("J8").Value = 1537
$ws.Range("J9").Value = 126
$ws.Range("J11").Value = 410
$ws.Range("J19").Value = 713
$ws.Range("J23").Value = 224
$ws.Range("J25").Value = 120
$ws.Range("J29").Value = 1327
$ws.Range("J30").Value = 88
$ws.Range("J31").Value = 230
$ws.Range("J33").Value = 1093
$ws.Range("J34").Value = 112
$ws.Range("J36").Value = 331
$ws.Range("J37").Value = 747
$ws.Range("J42").Value = 1047
$ws.Range("J43").Value = 207
$ws.Range("J49").Value = 155
$ws.Range("J51").Value = 300
$ws.Range("J53").Value = 353
$ws.Range("J54").Value = 459
$ws.Range("J60").Value = 142
$ws.Range("J63").Value = 80
$ws.Range("J65").Value = 601
$ws.Range("J66").Value = 73
$ws.Range("J67").Value = 916
$ws.Range("J72").Value = 96
$ws.Range("J76").Value = 361
$ws.Range("J78").Value = 288
$ws.Range("J79").Value = 685
$ws.Range("J80").Value = 41
$ws.Range("J83").Value = 482
$ws.Range("J84").Value = 203
$ws.Range("J85").Value = 1013
$ws.Range("J86").Value = 157
$ws.Range("J91").Value = 280
$ws.Range("J95").Value = 352
$ws.Range("J96").Value = 267
$ws.Range("J98").Value = 181
$ws.Range("J99").Value = 373
$ws.Range("J101").Value = 24317

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("J6").Value = 94
$ws.Range("J7").Value = 267

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("J3").Value = 217
$ws.Range("J6").Value = 228
$ws.Range("J7").Value = 708

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("J2").Value = 119
$ws.Range("J7").Value = 410

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("J2").Value = 269
$ws.Range("J3").Value = 359
$ws.Range("J6").Value = 293
$ws.Range("J7").Value = 1013

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("J3").Value = 44
$ws.Range("J7").Value = 353

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("J2").Value = 414
$ws.Range("J3").Value = 463
$ws.Range("J7").Value = 1537

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("J6").Value = 129
$ws.Range("J7").Value = 482

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("J2").Value = 258
$ws.Range("J3").Value = 361
$ws.Range("J4").Value = 47
$ws.Range("J7").Value = 1093

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("J2").Value = 120
$ws.Range("J4").Value = 17
$ws.Range("J7").Value = 352

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("J3").Value = 254
$ws.Range("J6").Value = 221
$ws.Range("J7").Value = 747

$ws = $wb.Worksheets.Item('New City')
$ws.Range("J3").Value = 170
$ws.Range("J6").Value = 216
$ws.Range("J7").Value = 601

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("J3").Value = 148
$ws.Range("J7").Value = 373

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("J3").Value = 33
$ws.Range("J7").Value = 88

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("J6").Value = 71
$ws.Range("J7").Value = 230

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("J3").Value = 341
$ws.Range("J6").Value = 254
$ws.Range("J7").Value = 916

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("J3").Value = 67
$ws.Range("J4").Value = 8
$ws.Range("J7").Value = 203

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("J6").Value = 89
$ws.Range("J7").Value = 155

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("J3").Value = 95
$ws.Range("J7").Value = 459

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("J2").Value = 404
$ws.Range("J3").Value = 467
$ws.Range("J6").Value = 335
$ws.Range("J7").Value = 1327

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("J3").Value = 208
$ws.Range("J6").Value = 278
$ws.Range("J7").Value = 713

$ws = $wb.Worksheets.Item('River North')
$ws.Range("J6").Value = 198
$ws.Range("J7").Value = 361

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("J2").Value = 56
$ws.Range("J6").Value = 69
$ws.Range("J7").Value = 185

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J2").Value = 223
$ws.Range("J3").Value = 203
$ws.Range("J7").Value = 1047

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("J6").Value = 86
$ws.Range("J7").Value = 288

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("J3").Value = 76
$ws.Range("J7").Value = 224

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("J3").Value = 117
$ws.Range("J4").Value = 10
$ws.Range("J7").Value = 280

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("J2").Value = 191
$ws.Range("J7").Value = 685

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("J3").Value = 105
$ws.Range("J7").Value = 331

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range("J6").Value = 42
$ws.Range("J7").Value = 112

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("J5").Value = 5
$ws.Range("J7").Value = 120

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("J2").Value = 34
$ws.Range("J7").Value = 181

$ws = $wb.Worksheets.Item('North Center')
$ws.Range("J6").Value = 44
$ws.Range("J7").Value = 73

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("J6").Value = 44
$ws.Range("J7").Value = 126

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("J2").Value = 56
$ws.Range("J7").Value = 194

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("J3").Value = 25
$ws.Range("J4").Value = 85
$ws.Range("J7").Value = 157

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("J3").Value = 80
$ws.Range("J7").Value = 300

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("J4").Value = 9
$ws.Range("J6").Value = 39
$ws.Range("J7").Value = 142

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("J3").Value = 41
$ws.Range("J6").Value = 123
$ws.Range("J7").Value = 207

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("J3").Value = 29
$ws.Range("J7").Value = 96

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range("J3").Value = 9
$ws.Range("J7").Value = 41

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range("J6").Value = 47
$ws.Range("J7").Value = 110
